$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 145
$ws1.Range("F5").Value = 45

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 145
$ws4.Range("F5").Value = 45
